$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "TEAM MEMBER" result cells for weeks 1-5 (G12, G19, G24, G29, G34) previously
# listed only one member ("Võ Hàn Trân Châu"); add the second member on a new line.
$newline = [char]10
$text = "Võ Hàn Trân Châu" + $newline + "Huỳnh Quốc Huy"

$targetCells = @("G12", "G19", "G24", "G29", "G34")
foreach ($ref in $targetCells) {
    $rng = $ws.Range($ref)
    $rng.Value2 = $text
    $rng.WrapText = $true
    $rng.EntireRow.RowHeight = 26.25
}

# Update the current selection to match the author's last interaction (range E9:E36).
$ws.Range("E9:E36").Select() | Out-Null
